$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from U1 into V1:X1, then set header labels
$ws.Range("U1").Copy()
$ws.Range("V1:X1").PasteSpecial(-4122)
$ws.Range("V1").Value = "poblacion2016"
$ws.Range("W1").Value = "pib"
$ws.Range("X1").Value = "quartile"

$ws.Range("V2").Value = 29305
$ws.Range("W2").Value = 32867.59
$ws.Range("X2").Value = "first quartile"
$ws.Range("V3").Value = 44434
$ws.Range("W3").Value = 88617.14
$ws.Range("X3").Value = "second quartile"
$ws.Range("V4").Value = 132041
$ws.Range("W4").Value = 288262.31
$ws.Range("X4").Value = "more third quartile"
$ws.Range("V5").Value = 103731
$ws.Range("W5").Value = 230106.32
$ws.Range("X5").Value = "more third quartile"
$ws.Range("V6").Value = 24862
$ws.Range("W6").Value = 43572.24
$ws.Range("X6").Value = "first quartile"
$ws.Range("V7").Value = 25448
$ws.Range("W7").Value = 32581.87
$ws.Range("X7").Value = "first quartile"
$ws.Range("V8").Value = 24302
$ws.Range("W8").Value = 240970.74
$ws.Range("X8").Value = "more third quartile"
$ws.Range("V9").Value = 74804
$ws.Range("W9").Value = 116284.87
$ws.Range("X9").Value = "third quartile"
$ws.Range("V10").Value = 19300
$ws.Range("W10").Value = 146767.18
$ws.Range("X10").Value = "third quartile"
$ws.Range("V11").Value = 253441
$ws.Range("W11").Value = 1849299.44
$ws.Range("X11").Value = "more third quartile"
$ws.Range("V12").Value = 92234
$ws.Range("W12").Value = 659769.89
$ws.Range("X12").Value = "more third quartile"
$ws.Range("V13").Value = 4589
$ws.Range("W13").Value = 13358.26
$ws.Range("X13").Value = "first quartile"
$ws.Range("V14").Value = 37929
$ws.Range("W14").Value = 52329.28
$ws.Range("X14").Value = "second quartile"
$ws.Range("V15").Value = 61193
$ws.Range("W15").Value = 128480.48
$ws.Range("X15").Value = "third quartile"
$ws.Range("V16").Value = 30546
$ws.Range("W16").Value = 42251.74
$ws.Range("X16").Value = "first quartile"
$ws.Range("V17").Value = 310582
$ws.Range("W17").Value = 1580563.51
$ws.Range("X17").Value = "more third quartile"
$ws.Range("V18").Value = 23342
$ws.Range("W18").Value = 42533.56
$ws.Range("X18").Value = "first quartile"
$ws.Range("V19").Value = 36470
$ws.Range("W19").Value = 56999.28
$ws.Range("X19").Value = "second quartile"
$ws.Range("V20").Value = 24139
$ws.Range("W20").Value = 47652.57
$ws.Range("X20").Value = "second quartile"
$ws.Range("V21").Value = 48920
$ws.Range("W21").Value = 62482.19
$ws.Range("X21").Value = "second quartile"
$ws.Range("V22").Value = 61553
$ws.Range("W22").Value = 205934.49
$ws.Range("X22").Value = "third quartile"
$ws.Range("V23").Value = 41524
$ws.Range("W23").Value = 100939.6
$ws.Range("X23").Value = "third quartile"
